$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5194237232208252
$ws.Range("B1").Value = 1.698515176773071
$ws.Range("C1").Value = 4.663849353790283
$ws.Range("D1").Value = 1.195181369781494
$ws.Range("E1").Value = 0.7456616759300232
